# "10Th - MB for single stock and added new group"
#
# The MarketBeat rank sheet keeps a rolling window of "last rated" columns,
# newest first (column B), with older days shifting right each time a new
# day is appended. This run layers in two more days of data - Jun_26
# (recorded as a pair of sibling columns) and Jun_27 - and appends two
# brand-new analyst rows (Benchmark, Evercore ISI) that just started
# coverage.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Make room: insert 3 new columns right after column A (the firm-name
#    column). This shifts the existing Jun_17/Jun_15/Jun_13/Jun_10 columns
#    (B:E) three places to the right (E:H), exactly like the sheet already
#    does every time a newer day is recorded.
# ---------------------------------------------------------------------
$ws.Range("B1:D1").EntireColumn.Insert()

# Keep the data columns at the sheet's usual 8-character width (matches
# the width already carried by the shifted-right columns).
$ws.Columns("C:H").ColumnWidth = 7.165

# ---------------------------------------------------------------------
# 2) Jun_26 is recorded first, as the pair of columns C and D.
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Row 6 - Goldman Sachs Group: raised target on 6/25/2018 (Jun_26 report).
$ws.Range("C6").Value = "6/25/2018,Raises Target,Buy,`$25.00 -> `$28.00"
$ws.Range("D6").Value = "6/25/2018,Raises Target,Buy,`$25.00 -> `$28.00"
$ws.Range("D6").Interior.ColorIndex = 35

# Row 21 - Wells Fargo & Co: upgraded on 6/18/2018 (reported with Jun_26).
$ws.Range("C21").Value = "6/18/2018,Upgrades,Underperform -> Market Perform,"
$ws.Range("D21").Value = "6/18/2018,Upgrades,Underperform -> Market Perform,"

# ---------------------------------------------------------------------
# 3) Two brand-new analyst rows that just started coverage - only the
#    firm name plus the three newest columns are populated.
# ---------------------------------------------------------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# ---------------------------------------------------------------------
# 4) Jun_27 is recorded second, as column B.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Jun_27"

for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
}

# Row 19 - Citigroup: raised target on 6/27/2018 (Jun_27 report).
$ws.Range("B19").Value = "6/27/2018,Raises Target,Buy,`$25.00 -> `$27.00"
